$d = $word.ActiveDocument

# Commit "paragrafKORIST": add a new "KORIST" (benefits) section, right
# after the "sintaksa" section and before the document's closing
# (bookmark-only) paragraph.
#
# The whole thing is guarded so the script is idempotent: if the
# "KORIST" heading is already present in the document nothing is
# changed.

$alreadyPresent = $d.Content.Find.Execute("KORIST", $true, $true, $false,
    $false, $false, $true, 1, $false, "", 0)

if (-not $alreadyPresent) {

    $anchorText = "zamjena vrijednosti unutar varijabli " +
        "što možemo učiniti jednostavnim izrazom (a,b = b,a)."

    # Preferred insertion point: right after the last sentence of the
    # "sintaksa" section.
    $rng = $d.Content
    $found = $rng.Find.Execute($anchorText, $true, $true, $false, $false,
        $false, $true, 1, $false, "", 0)

    if (-not $found) {
        # Fallback: anchor on the paragraph right before the document's
        # last paragraph (the trailing bookmark-only paragraph at the
        # very end), so the new section still lands at the end of the
        # article body.
        $paraCount = $d.Paragraphs.Count
        $rng = $d.Paragraphs.Item($paraCount - 1).Range
        $rng.Collapse(0)
    } else {
        $rng.Collapse(0)
    }

    $rng.InsertParagraphAfter()
    $rng.Start = $rng.Start + 1
    $rng.End = $rng.Start

    $bodyLines = @(
        "Od 2003 Python spada u top 10 najpopularnijih  programskih jezika. 2017. je dobio status trečeg najpopularnijeg jezika koji nema sintaksu C-a.",
        "Python može služiti kao skripterski jezik za web aplikacijekoristeči web frameworkove kao što je Django.",
        "Library poput NumPy, SciPy i Matplotlib omogučuju korištenje pythona u znanstvenoj obradi podataka.",
        "Mnogi operativni sustavi dodaju Python kao standardnu komponentu.",
        "LibreOffice će uskoro zamjeniti svoj kod u Javi sa Pythonom. "
    )

    $allLines = New-Object System.Collections.ArrayList
    [void]$allLines.Add("KORIST")
    foreach ($l in $bodyLines) { [void]$allLines.Add($l) }
    $blob = [string]::Join([char]13, $allLines)

    $rng.InsertAfter($blob)

    # The first of the newly-inserted paragraphs is the "KORIST" heading -
    # give it the same style used by the other section headings.
    $headingPara = $rng.Paragraphs(1)
    $headingPara.Range.Style = "Heading2"
}
